$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), copying the format of the
# existing header style (bold, centered, bordered) from H1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in I and J numeric values for each data row (2-51)
$iValues = @(8,10,9,6,8,7,8,9,7,4,5,9,7,8,5,8,8,9,6,8,6,1,7,9,8,9,7,9,9,10,7,7,2,7,8,7,6,7,10,8,6,6,9,6,8,6,5,6,5,3)
$jValues = @(8,10,9,7,8,7,8,9,8,6,7,9,8,9,5,8,9,9,7,8,7,1,7,9,8,9,8,9,9,11,7,8,4,7,8,7,7,8,11,9,8,7,9,8,8,9,8,7,6,4)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}

